$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111863040
$ws.Range("B2").Value = 90687
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 5964
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "Fjällig taggsvamp s.str."
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "Sarcodon imbricatus s.str."
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "(L.:Fr.) P.Karst."
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = $null
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = $null
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q2").Value = 655235.4020021557
$ws.Range("R2").Value = 6634878.090185729
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "10:49"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "10:49"
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = "Halv häxring, 3 m i diameter"

# Row 3
$ws.Range("A3").Value = 111863073
$ws.Range("B3").Value = 88899
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 3286
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "Flattoppad klubbsvamp"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "Clavariadelphus truncatus"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "(Quél.) Donk"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "2"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("P3").NumberFormat = "@"
$ws.Range("P3").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q3").Value = 655228.290648401
$ws.Range("R3").Value = 6634879.303300899
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "10:50"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "10:50"
$ws.Range("AC3").NumberFormat = "@"
$ws.Range("AC3").Value = $null

# Row 4
$ws.Range("A4").Value = 111863288
$ws.Range("B4").Value = 85062
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 249278
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "Barrviolspindling"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Cortinarius harcynicus"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "(Pers.) M.M.Moser"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "fruktkroppar"
$ws.Range("P4").NumberFormat = "@"
$ws.Range("P4").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q4").Value = 655134.5683182024
$ws.Range("R4").Value = 6634792.815828164
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "11:02"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "11:02"
$ws.Range("AC4").NumberFormat = "@"
$ws.Range("AC4").Value = $null

# Row 5
$ws.Range("A5").Value = 111863045
$ws.Range("B5").Value = 88899
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 3286
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "Flattoppad klubbsvamp"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Clavariadelphus truncatus"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "(Quél.) Donk"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "11"
$ws.Range("J5").NumberFormat = "@"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("P5").NumberFormat = "@"
$ws.Range("P5").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q5").Value = 655233.932825509
$ws.Range("R5").Value = 6634889.105171775
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "10:50"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "10:50"
$ws.Range("AC5").NumberFormat = "@"
$ws.Range("AC5").Value = "Under gran i svacka"

# Row 6
$ws.Range("A6").Value = 111863218
$ws.Range("B6").Value = 90021
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 6031
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "Blomkålssvamp"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Sparassis crispa"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "(Wulfen:Fr.) Fr."
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("P6").NumberFormat = "@"
$ws.Range("P6").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q6").Value = 655137.9235184891
$ws.Range("R6").Value = 6634821.151011234
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "10:53"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "10:53"
$ws.Range("AC6").NumberFormat = "@"
$ws.Range("AC6").Value = $null

# Row 7
$ws.Range("A7").Value = 111863269
$ws.Range("B7").Value = 85062
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 249278
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "Barrviolspindling"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "Cortinarius harcynicus"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "(Pers.) M.M.Moser"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "4"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "fruktkroppar"
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q7").Value = 655135.2812587479
$ws.Range("R7").Value = 6634799.89438487
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = "11:02"
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = "11:02"
$ws.Range("AC7").NumberFormat = "@"
$ws.Range("AC7").Value = "4 ex i gräsglänta under gran och tall."

# Row 8
$ws.Range("A8").Value = 111863001
$ws.Range("B8").Value = 90332
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4769
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Svavelriska"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Lactarius scrobiculatus"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "(Scop.:Fr.) Fr."
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "1"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "fruktkroppar"
$ws.Range("P8").NumberFormat = "@"
$ws.Range("P8").Value = "Charlottenberg, Upl"
$ws.Range("Q8").Value = 655217.6931657954
$ws.Range("R8").Value = 6634939.780080916
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "10:47"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "10:47"
$ws.Range("AC8").NumberFormat = "@"
$ws.Range("AC8").Value = $null

# Row 9
$ws.Range("A9").Value = 111863402
$ws.Range("B9").Value = 90687
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 5964
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Fjällig taggsvamp s.str."
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Sarcodon imbricatus s.str."
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(L.:Fr.) P.Karst."
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "1"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "fruktkroppar"
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "Charlottenberg (Charlottenberg), Upl"
$ws.Range("Q9").Value = 655199.5794486763
$ws.Range("R9").Value = 6634769.85474884
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "11:02"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "11:02"
$ws.Range("AC9").NumberFormat = "@"
$ws.Range("AC9").Value = $null
